# Update res_bus/vm_pu.xlsx values ("case with 380 kV done")
# Slack bus voltage setpoint changed from 1.05 to 1.02 pu, and the
# dependent bus voltage magnitudes (columns B-F, I-N, rows 2-25) were
# recomputed accordingly. Columns A, G, H are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.045152066849036
$ws.Range("D2").Value = 1.052526514471977
$ws.Range("E2").Value = 1.05259724043079
$ws.Range("F2").Value = 1.062595122994376
$ws.Range("I2").Value = 1.037586549292694
$ws.Range("J2").Value = 1.050213330853448
$ws.Range("K2").Value = 1.055275002408555
$ws.Range("L2").Value = 1.055345532738893
$ws.Range("M2").Value = 1.065316044784235
$ws.Range("N2").Value = 1.051704754244413

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.046153867393417
$ws.Range("D3").Value = 1.053427385776767
$ws.Range("E3").Value = 1.053493831880952
$ws.Range("F3").Value = 1.063590751221138
$ws.Range("I3").Value = 1.037737667355049
$ws.Range("J3").Value = 1.050862418683493
$ws.Range("K3").Value = 1.055988511132466
$ws.Range("L3").Value = 1.056054786742762
$ws.Range("M3").Value = 1.066126065192735
$ws.Range("N3").Value = 1.052354763853628

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046802018643817
$ws.Range("D4").Value = 1.054010579883392
$ws.Range("E4").Value = 1.054074339351981
$ws.Range("F4").Value = 1.064235513230369
$ws.Range("I4").Value = 1.037833253474823
$ws.Range("J4").Value = 1.051281733084722
$ws.Range("K4").Value = 1.056449834989849
$ws.Range("L4").Value = 1.056513439100409
$ws.Range("M4").Value = 1.066650109526314
$ws.Range("N4").Value = 1.052774673729403

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.047074481523156
$ws.Range("D5").Value = 1.054255818328413
$ws.Range("E5").Value = 1.054318468131277
$ws.Range("F5").Value = 1.064506695604809
$ws.Range("I5").Value = 1.037872911125407
$ws.Range("J5").Value = 1.051457847143459
$ws.Range("K5").Value = 1.05664368754388
$ws.Range("L5").Value = 1.056706188094144
$ws.Range("M5").Value = 1.066870394579837
$ws.Range("N5").Value = 1.052951037890303

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047120228054956
$ws.Range("D6").Value = 1.054296998644443
$ws.Range("E6").Value = 1.054359463290202
$ws.Range("F6").Value = 1.064552235585521
$ws.Range("I6").Value = 1.03787953890982
$ws.Range("J6").Value = 1.051487407727558
$ws.Range("K6").Value = 1.056676231055626
$ws.Range("L6").Value = 1.056738547459162
$ws.Range("M6").Value = 1.06690738002345
$ws.Range("N6").Value = 1.052980640453822

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046805659385114
$ws.Range("D7").Value = 1.05401385652273
$ws.Range("E7").Value = 1.054077601086135
$ws.Range("F7").Value = 1.064239136294907
$ws.Range("I7").Value = 1.03783378545297
$ws.Range("J7").Value = 1.051284086981575
$ws.Range("K7").Value = 1.056452425604014
$ws.Range("L7").Value = 1.056516014891967
$ws.Range("M7").Value = 1.066653053078238
$ws.Range("N7").Value = 1.052777030969059

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045490646388903
$ws.Range("D8").Value = 1.052830911599123
$ws.Range("E8").Value = 1.052900174029935
$ws.Range("F8").Value = 1.062931491116278
$ws.Range("I8").Value = 1.037638074912416
$ws.Range("J8").Value = 1.050432835267461
$ws.Range("K8").Value = 1.055516210879794
$ws.Range("L8").Value = 1.055585286417432
$ws.Range("M8").Value = 1.065589813691777
$ws.Range("N8").Value = 1.051924570379871

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043172832748289
$ws.Range("D9").Value = 1.050748525152882
$ws.Range("E9").Value = 1.050828144601189
$ws.Range("F9").Value = 1.060631310371556
$ws.Range("I9").Value = 1.037276408979596
$ws.Range("J9").Value = 1.048927580454169
$ws.Range("K9").Value = 1.053863730217319
$ws.Range("L9").Value = 1.053943097872803
$ws.Range("M9").Value = 1.063715577137104
$ws.Range("N9").Value = 1.050417177932089

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041627256214592
$ws.Range("D10").Value = 1.04936174681768
$ws.Range("E10").Value = 1.049448697607285
$ws.Range("F10").Value = 1.059100641543621
$ws.Range("I10").Value = 1.037024043230302
$ws.Range("J10").Value = 1.047920599066028
$ws.Range("K10").Value = 1.052760275028787
$ws.Range("L10").Value = 1.052846923833317
$ws.Range("M10").Value = 1.062465684991091
$ws.Range("N10").Value = 1.049408766514874

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040957922959934
$ws.Range("D11").Value = 1.048761618959317
$ws.Range("E11").Value = 1.048851845633078
$ws.Range("F11").Value = 1.058438516951805
$ws.Range("I11").Value = 1.036912104430586
$ws.Range("J11").Value = 1.047483748933222
$ws.Range("K11").Value = 1.052282048969754
$ws.Range("L11").Value = 1.052371950119412
$ws.Range("M11").Value = 1.061924382910862
$ws.Range("N11").Value = 1.048971296004773

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.0407092898445
$ws.Range("D12").Value = 1.048538759147472
$ws.Range("E12").Value = 1.048630217873804
$ws.Range("F12").Value = 1.058192675161426
$ws.Range("I12").Value = 1.036870125900114
$ws.Range("J12").Value = 1.047321360705711
$ws.Range("K12").Value = 1.052104351707082
$ws.Range("L12").Value = 1.052195475880989
$ws.Range("M12").Value = 1.06172330635887
$ws.Range("N12").Value = 1.048808677167352

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040762623072028
$ws.Range("D13").Value = 1.048586560858335
$ws.Range("E13").Value = 1.048677754615181
$ws.Range("F13").Value = 1.058245404486812
$ws.Range("I13").Value = 1.036879148505411
$ws.Range("J13").Value = 1.047356199071151
$ws.Range("K13").Value = 1.052142471204448
$ws.Range("L13").Value = 1.052233332358585
$ws.Range("M13").Value = 1.061766438531826
$ws.Range("N13").Value = 1.048843565007267

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040937371131121
$ws.Range("D14").Value = 1.04874319618313
$ws.Range("E14").Value = 1.048833524387697
$ws.Range("F14").Value = 1.058418193532338
$ws.Range("I14").Value = 1.036908642618065
$ws.Range("J14").Value = 1.047470328369768
$ws.Range("K14").Value = 1.052267361735945
$ws.Range("L14").Value = 1.052357363676016
$ws.Range("M14").Value = 1.061907762118369
$ws.Range("N14").Value = 1.048957856382579

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041045037472835
$ws.Range("D15").Value = 1.04883971160998
$ws.Range("E15").Value = 1.048929508539546
$ws.Range("F15").Value = 1.058524667943235
$ws.Range("I15").Value = 1.036926761996855
$ws.Range("J15").Value = 1.047540630965765
$ws.Range("K15").Value = 1.052344302611764
$ws.Range("L15").Value = 1.052433777143859
$ws.Range("M15").Value = 1.06199483456263
$ws.Range("N15").Value = 1.049028258816326

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041671675824702
$ws.Range("D16").Value = 1.049401582908141
$ws.Range("E16").Value = 1.04948831842497
$ws.Range("F16").Value = 1.059144598657004
$ws.Range("I16").Value = 1.037031416187587
$ws.Range("J16").Value = 1.047949574137475
$ws.Range("K16").Value = 1.052792004471134
$ws.Range("L16").Value = 1.052878439525892
$ws.Range("M16").Value = 1.062501607583633
$ws.Range("N16").Value = 1.049437782734247

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042064726078323
$ws.Range("D17").Value = 1.049754125936118
$ws.Range("E17").Value = 1.049838968257977
$ws.Range("F17").Value = 1.059533643479043
$ws.Range("I17").Value = 1.03709635050688
$ws.Range("J17").Value = 1.048205874155976
$ws.Range("K17").Value = 1.053072723346479
$ws.Range("L17").Value = 1.05315727847276
$ws.Range("M17").Value = 1.06281946920007
$ws.Range("N17").Value = 1.049694446728172

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042293977087957
$ws.Range("D18").Value = 1.049959792641279
$ws.Range("E18").Value = 1.050043540473815
$ws.Range("F18").Value = 1.059760630859423
$ws.Range("I18").Value = 1.03713396864052
$ws.Range("J18").Value = 1.04835529040716
$ws.Range("K18").Value = 1.053236421088565
$ws.Range("L18").Value = 1.053319889290356
$ws.Range("M18").Value = 1.06300486371858
$ws.Range("N18").Value = 1.049844075167571

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042372144296595
$ws.Range("D19").Value = 1.050029925467962
$ws.Range("E19").Value = 1.050113301766222
$ws.Range("F19").Value = 1.059838038540118
$ws.Range("I19").Value = 1.037146751850121
$ws.Range("J19").Value = 1.04840622403199
$ws.Range("K19").Value = 1.053292230841907
$ws.Range("L19").Value = 1.053375330059858
$ws.Range("M19").Value = 1.063068076945618
$ws.Range("N19").Value = 1.04989508112399

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042022556388028
$ws.Range("D20").Value = 1.049716297849286
$ws.Range("E20").Value = 1.049801342281432
$ws.Range("F20").Value = 1.059491896013359
$ws.Range("I20").Value = 1.03708941024293
$ws.Range("J20").Value = 1.048178383778282
$ws.Range("K20").Value = 1.053042609085851
$ws.Range("L20").Value = 1.053127364917072
$ws.Range("M20").Value = 1.062785366560003
$ws.Range("N20").Value = 1.049666917310988

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040885912537765
$ws.Range("D21").Value = 1.048697069475532
$ws.Range("E21").Value = 1.048787652142764
$ws.Range("F21").Value = 1.058367308676799
$ws.Range("I21").Value = 1.03689996835994
$ws.Range("J21").Value = 1.04743672351621
$ws.Range("K21").Value = 1.052230586315614
$ws.Range("L21").Value = 1.052320840859878
$ws.Range("M21").Value = 1.061866146193533
$ws.Range("N21").Value = 1.048924203806275

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040171184719111
$ws.Range("D22").Value = 1.048056555281234
$ws.Range("E22").Value = 1.048150708257225
$ws.Range("F22").Value = 1.057660819383907
$ws.Range("I22").Value = 1.036778547761207
$ws.Range("J22").Value = 1.046969703042295
$ws.Range("K22").Value = 1.051719672301163
$ws.Range("L22").Value = 1.051813470650158
$ws.Range("M22").Value = 1.061288121836518
$ws.Range("N22").Value = 1.048456520109719

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040550082428729
$ws.Range("D23").Value = 1.048396073816844
$ws.Range("E23").Value = 1.048488325739344
$ws.Range("F23").Value = 1.058035287147492
$ws.Range("I23").Value = 1.036843133981578
$ws.Range("J23").Value = 1.047217346446883
$ws.Range("K23").Value = 1.051990551672806
$ws.Range("L23").Value = 1.052082463180008
$ws.Range("M23").Value = 1.06159455025541
$ws.Range("N23").Value = 1.048704515196346

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042041611082233
$ws.Range("D24").Value = 1.049733390627402
$ws.Range("E24").Value = 1.049818343704398
$ws.Range("F24").Value = 1.059510759697994
$ws.Range("I24").Value = 1.037092547043436
$ws.Range("J24").Value = 1.048190805741397
$ws.Range("K24").Value = 1.053056216550125
$ws.Range("L24").Value = 1.053140881661745
$ws.Range("M24").Value = 1.062800776102259
$ws.Range("N24").Value = 1.049679356914715

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043772110145931
$ws.Range("D25").Value = 1.051286615456694
$ws.Range("E25").Value = 1.051363482204071
$ws.Range("F25").Value = 1.061225475512897
$ws.Range("I25").Value = 1.037371894944735
$ws.Range("J25").Value = 1.049317341238249
$ws.Range("K25").Value = 1.054291257073635
$ws.Range("L25").Value = 1.054367889972204
$ws.Range("M25").Value = 1.064200186636891
$ws.Range("N25").Value = 1.050807492221189
